$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Date) holds text values like "2025-12-24" that Excel would
# otherwise auto-convert into a date serial number. Force the cells to
# Text format first, then restore the default "Normal" style afterwards
# so the final cells carry no explicit style (matching the source rows).
$ws.Range("B2:B4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Thai League 1'
$ws.Range("B2").Value = '2025-12-24'
$ws.Range("C2").Value = '08:00:00'
$ws.Range("D2").Value = 'Ratchaburi'
$ws.Range("E2").Value = 'Chiangrai Utd'
$ws.Range("F2").Value = 1.51
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 8.199999999999999
$ws.Range("J2").Value = 3.85
$ws.Range("K2").Value = 5.7
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.01
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.72
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.01
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.14
$ws.Range("W2").Value = 2.36
$ws.Range("X2").Value = 28
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("A3").Value = 'Thai League 1'
$ws.Range("B3").Value = '2025-12-24'
$ws.Range("C3").Value = '08:00:00'
$ws.Range("D3").Value = 'Sukhothai'
$ws.Range("E3").Value = 'Buriram Utd'
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 11.5
$ws.Range("H3").Value = 1.39
$ws.Range("I3").Value = 1.48
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 5.8
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 2.36
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 2.36
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 2.46
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 1.09
$ws.Range("X3").Value = 34
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 13.5
$ws.Range("AA3").Value = 18
$ws.Range("AB3").Value = 46
$ws.Range("AC3").Value = 18
$ws.Range("AD3").Value = 15
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 48
$ws.Range("AH3").Value = 34
$ws.Range("AI3").Value = 48
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("A4").Value = 'Thai League 1'
$ws.Range("B4").Value = '2025-12-24'
$ws.Range("C4").Value = '09:00:00'
$ws.Range("D4").Value = 'BG Pathumthani United'
$ws.Range("E4").Value = 'Dragon Pathumwan Kanchana'
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Value = 950
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 1.24
$ws.Range("Q4").Value = 1.01
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0
$ws.Range("AO4").Value = 0

# Restore default styling on the date column (removes the temporary
# Text number format applied above, back to General/no explicit style).
$ws.Range("B2:B4").Style = "Normal"

